$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

# Copy formatting (bold font, borders, alignment) from the cell above in column A
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = "Fullness"

for ($col = 2; $col -le 62; $col++) {
    $ws.Cells.Item($row, $col).Value = 0
}
